# Update "想去人数" (want-to-go count) values on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 2095
$ws1.Range("F5").Value = 370
$ws1.Range("F6").Value = 651
$ws1.Range("F8").Value = 2083
$ws1.Range("F9").Value = 10768
$ws1.Range("F10").Value = 181
$ws1.Range("F11").Value = 162
$ws1.Range("F12").Value = 290
$ws1.Range("F13").Value = 205
$ws1.Range("F15").Value = 8989
$ws1.Range("F16").Value = 1116
$ws1.Range("F17").Value = 729
$ws1.Range("F18").Value = 5275
$ws1.Range("F19").Value = 72
$ws1.Range("F20").Value = 3353

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 2095
$ws4.Range("F5").Value = 370
$ws4.Range("F6").Value = 651
$ws4.Range("F9").Value = 2083
$ws4.Range("F12").Value = 10768
$ws4.Range("F13").Value = 181
$ws4.Range("F14").Value = 162
$ws4.Range("F15").Value = 290
$ws4.Range("F16").Value = 205
$ws4.Range("F18").Value = 8989
$ws4.Range("F19").Value = 1116
$ws4.Range("F20").Value = 729
$ws4.Range("F21").Value = 5275
$ws4.Range("F22").Value = 72
$ws4.Range("F23").Value = 3353
